# Apply cryptos.xlsx data refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.858.85"
$ws.Range("E2").Value = "  -4.73%  "
$ws.Range("D3").Value = "2.209.35"
$ws.Range("E3").Value = "  -6.24%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.84"
$ws.Range("E6").Value = "  -7.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.588"
$ws.Range("E7").Value = "  -6.68%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -8.15%  "
$ws.Range("E10").Value = "  -9.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.09"
$ws.Range("E11").Value = "  -2.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0823"
$ws.Range("E12").Value = "  -9.92%  "
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").Value = "2.544.12"
$ws.Range("E15").Value = "  -6.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.856"
$ws.Range("E16").Value = "  -11.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.18"
$ws.Range("E17").Value = "  -6.59%  "
$ws.Range("D18").Value = "2.207.28"
$ws.Range("E18").Value = "  -6.55%  "
$ws.Range("D19").Value = "42.719.29"
$ws.Range("E19").Value = "  -5.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.32"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").Value = "0.0₃0959"
$ws.Range("E21").Value = "  -9.49%  "
$ws.Range("E22").Value = "  -10.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.13"
$ws.Range("E23").Value = "  -10.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.13"
$ws.Range("E24").Value = "  -10.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "235.56"
$ws.Range("E25").Value = "  -8.92%  "
$ws.Range("E26").Value = "  -7.56%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.95"
$ws.Range("E30").Value = "  -9.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.27"
$ws.Range("E31").Value = "  -12.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0884"
$ws.Range("E32").Value = "  -9.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.49"
$ws.Range("E33").Value = "  -7.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.16"
$ws.Range("E34").Value = "  -7.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "153.84"
$ws.Range("E35").Value = "  -8.03%  "
$ws.Range("E36").Value = "  -7.11%  "
$ws.Range("E37").Value = "  +6.64%  "
$ws.Range("E38").Value = "  -6.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.91"
$ws.Range("E39").Value = "  +7.49%  "
$ws.Range("E40").Value = "  -7.90%  "
$ws.Range("E41").Value = "  -5.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.75"
$ws.Range("E42").Value = "  -4.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0323"
$ws.Range("E43").Value = "  -7.92%  "
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "1.822.38"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.15"
$ws.Range("E46").Value = "  -4.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.21"
$ws.Range("E47").Value = "  -11.90%  "
$ws.Range("E48").Value = "  -9.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.34"
$ws.Range("E49").Value = "  -5.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "60.60"
$ws.Range("E50").Value = "  -12.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.21"
$ws.Range("E51").Value = "  -9.54%  "
